# Apply updated crypto price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.445.16"
$ws.Range("D3").Value = "3.691.69"
$ws.Range("E3").Value = "  -3.10%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "686.90"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.20%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "161.96"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -5.56%  "
$ws.Range("D7").Value = "3.690.23"
$ws.Range("E7").Value = "  -3.14%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -5.51%  "
$ws.Range("E10").Value = "  -8.39%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "7.38"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.63%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.440"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -7.78%  "
$ws.Range("E13").Value = "  -5.91%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.311.53"
$ws.Range("E14").Value = "  -3.11%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "33.10"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -8.15%  "
$ws.Range("D16").Value = "3.690.25"
$ws.Range("E16").Value = "  -3.27%  "
$ws.Range("D17").Value = "69.490.40"
$ws.Range("E17").Value = "  -2.55%  "
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("E19").Value = "  -7.95%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "6.54"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -9.01%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "476.88"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -7.26%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "9.96"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -5.19%  "
$ws.Range("E23").Value = "  -7.63%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "79.94"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -4.92%  "
$ws.Range("D25").Value = "3.833.06"
$ws.Range("E25").Value = "  -3.13%  "
$ws.Range("E26").Value = "  -9.46%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  -7.99%  "
$ws.Range("E29").Value = "  -9.79%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.81"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -10.90%  "
$ws.Range("E31").Value = "  -10.11%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "6.80"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -7.28%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.06"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -7.73%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "26.95"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -7.38%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.166"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -4.20%  "
$ws.Range("D37").Value = "3.654.80"
$ws.Range("E37").Value = "  -3.07%  "
$ws.Range("E38").Value = "  -8.47%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "6.25"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -3.81%  "
$ws.Range("E40").Value = "  -3.20%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0921"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -8.85%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  -0.03%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.954"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -6.27%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "163.65"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -5.66%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "48.28"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.54%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "29.80"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("E48").Value = "  -15.35%  "
$ws.Range("E49").Value = "  -3.92%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.000282"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -8.83%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.12"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.09%  "
